$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the refreshed crypto data.
# Leading apostrophe forces text interpretation so Excel doesn't coerce
# strings like "1.002" or "5.890" into numbers (which would drop the
# trailing zero / change formatting), matching the original inline-string cells.

$ws.Range("D2").Formula = "'27.473.53"
$ws.Range("E2").Formula = "'  +0.24%  "
$ws.Range("D3").Formula = "'1.736.76"
$ws.Range("E3").Formula = "'  -0.73%  "
$ws.Range("D4").Formula = "'1.002"
$ws.Range("E4").Formula = "'  +0.02%  "
$ws.Range("D5").Formula = "'323.06"
$ws.Range("E5").Formula = "'  +0.61%  "
$ws.Range("E6").Formula = "'  +0.06%  "
$ws.Range("D7").Formula = "'0.4522"
$ws.Range("E7").Formula = "'  +6.99%  "
$ws.Range("D8").Formula = "'0.3528"
$ws.Range("E8").Formula = "'  -1.89%  "
$ws.Range("D9").Formula = "'0.07379"
$ws.Range("E9").Formula = "'  -2.05%  "
$ws.Range("D11").Formula = "'1.071"
$ws.Range("E11").Formula = "'  -2.00%  "
$ws.Range("D12").Formula = "'1.001"
$ws.Range("E12").Formula = "'  -0.05%  "
$ws.Range("E13").Formula = "'  -1.70%  "
$ws.Range("D14").Formula = "'5.890"
$ws.Range("E14").Formula = "'  -2.33%  "
$ws.Range("D15").Formula = "'7.027"
$ws.Range("E15").Formula = "'  -2.48%  "
$ws.Range("D16").Formula = "'1.740.78"
$ws.Range("E16").Formula = "'  -0.47%  "
$ws.Range("D17").Formula = "'91.35"
$ws.Range("E17").Formula = "'  -0.02%  "
$ws.Range("D18").Formula = "'0.00001052"
$ws.Range("E18").Formula = "'  -1.54%  "
$ws.Range("D19").Formula = "'0.06333"
$ws.Range("E19").Formula = "'  -0.28%  "
$ws.Range("E20").Formula = "'  -0.02%  "
$ws.Range("D21").Formula = "'16.52"
$ws.Range("E21").Formula = "'  -2.96%  "
$ws.Range("D22").Formula = "'5.708"
$ws.Range("E22").Formula = "'  -2.89%  "
$ws.Range("D23").Formula = "'27.504.42"
$ws.Range("E23").Formula = "'  +0.22%  "
$ws.Range("E24").Formula = "'  -1.04%  "
$ws.Range("D25").Formula = "'2.084"
$ws.Range("E25").Formula = "'  -0.30%  "
$ws.Range("D26").Formula = "'161.67"
$ws.Range("E26").Formula = "'  +0.34%  "
$ws.Range("D27").Formula = "'19.94"
$ws.Range("E27").Formula = "'  -1.53%  "
$ws.Range("D28").Formula = "'1.938.57"
$ws.Range("E28").Formula = "'  -0.61%  "
$ws.Range("D29").Formula = "'124.29"
$ws.Range("E29").Formula = "'  +0.93%  "
$ws.Range("E30").Formula = "'  -4.95%  "
$ws.Range("D31").Formula = "'1.042"
$ws.Range("D32").Formula = "'0.09053"
$ws.Range("E32").Formula = "'  +2.47%  "
$ws.Range("D33").Formula = "'3.648"
$ws.Range("E33").Formula = "'  +0.17%  "
$ws.Range("D34").Formula = "'5.360"
$ws.Range("E34").Formula = "'  -3.26%  "
$ws.Range("E35").Formula = "'  -0.62%  "
$ws.Range("D36").Formula = "'11.55"
$ws.Range("E36").Formula = "'  -5.44%  "
$ws.Range("E37").Formula = "'  -1.06%  "
$ws.Range("D38").Formula = "'0.2051"
$ws.Range("E38").Formula = "'  -2.26%  "
$ws.Range("D39").Formula = "'0.6211"
$ws.Range("E39").Formula = "'  -1.71%  "
$ws.Range("D40").Formula = "'4.867"
$ws.Range("E40").Formula = "'  -1.17%  "
$ws.Range("D41").Formula = "'1.187"
$ws.Range("E41").Formula = "'  +0.57%  "
$ws.Range("D42").Formula = "'1.372"
$ws.Range("D43").Formula = "'7.667"
$ws.Range("E43").Formula = "'  -2.37%  "
$ws.Range("D44").Formula = "'13.00"
$ws.Range("E44").Formula = "'  -3.19%  "
$ws.Range("E45").Formula = "'  +0.41%  "
$ws.Range("D46").Formula = "'0.5767"
$ws.Range("E46").Formula = "'  -1.42%  "
$ws.Range("D47").Formula = "'122.06"
$ws.Range("E47").Formula = "'  -0.29%  "
$ws.Range("D48").Formula = "'1.925"
$ws.Range("E48").Formula = "'  -2.15%  "
$ws.Range("D49").Formula = "'0.06841"
$ws.Range("E49").Formula = "'  +0.60%  "
$ws.Range("E50").Formula = "'  -4.10%  "
$ws.Range("E51").Formula = "'  -2.85%  "
